# Daily attendance processing - 2026-01-08 16:06:13
#
# The "Recorded By" column (G) lists "System, <email>" for every session
# that was auto-recorded together with a human editor. Flip the display
# order to "<email>, System" throughout the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldText = "System, dnasr281@gmail.com"
$newText = "dnasr281@gmail.com, System"

# Column G ("Recorded By") holds the used range for this sheet; restrict
# the replace to that column so only matching "Recorded By" cells are
# touched (xlWhole => only cells whose entire content equals $oldText).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row
$colG = $ws.Range("G1:G$lastRow")

$colG.Replace($oldText, $newText, -4163, 1, $false, $false, $true, $true) | Out-Null
